# Insert a new row of data (LeetCode problem 183 "Customers Who Never Order")
# at row 5, pushing the existing "row 5" formatting content down visually by
# filling in the previously-empty cells on that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use the row above (row 4) as a formatting template, then overwrite values.
$ws.Range("A4:H4").Copy()
$ws.Range("A5:H5").PasteSpecial()

$ws.Range("A5").Value = 183
$ws.Range("B5").Value = "Customers Who Never Order"
$ws.Range("C5").Value = "SELECT"
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = "Easy"
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = "✅"
$ws.Range("H5").Value = "Given 1 sol and didn’t see solutions"

$ws.Rows.Item(5).RowHeight = 30

$ws.Range("E11").Select()
